# Update TPM-derived NATMI ligand-receptor metrics on the active sheet
# (Tgfb1-Acvrl1) to reflect the new TPM expression values. Only numeric
# cells affected by the new TPM calculation are updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 29.80827733333334
$ws.Range("N2").Value = 89.42483200000001
$ws.Range("O2").Value = 0.4866975737940222
$ws.Range("P2").Value = 0.4866975737940221
$ws.Range("Q2").Value = 1070.325188234176
$ws.Range("R2").Value = 9632.926694107584
$ws.Range("S2").Value = 0.2972304497777762
$ws.Range("T2").Value = 0.2972304497777761
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.437868100938039
$ws.Range("P3").Value = 0.437868100938039
$ws.Range("Q3").Value = 962.9414297359789
$ws.Range("R3").Value = 8666.47286762381
$ws.Range("S3").Value = 0.2674098651665656
$ws.Range("T3").Value = 0.2674098651665655
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("M4").Value = 4.620050333333332
$ws.Range("O4").Value = 0.07543432526793886
$ws.Range("P4").Value = 0.07543432526793886
$ws.Range("Q4").Value = 165.892050297943
$ws.Range("R4").Value = 1493.028452681487
$ws.Range("S4").Value = 0.04606839983459956
$ws.Range("T4").Value = 0.04606839983459955
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 29.80827733333334
$ws.Range("N5").Value = 89.42483200000001
$ws.Range("O5").Value = 0.4866975737940222
$ws.Range("P5").Value = 0.4866975737940221
$ws.Range("Q5").Value = 508.1900725993531
$ws.Range("R5").Value = 4573.710653394177
$ws.Range("S5").Value = 0.1411249268089338
$ws.Range("T5").Value = 0.1411249268089338
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.437868100938039
$ws.Range("P6").Value = 0.437868100938039
$ws.Range("S6").Value = 0.1269661223398657
$ws.Range("T6").Value = 0.1269661223398656
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("M7").Value = 4.620050333333332
$ws.Range("O7").Value = 0.07543432526793886
$ws.Range("P7").Value = 0.07543432526793886
$ws.Range("Q7").Value = 78.76549483400755
$ws.Range("R7").Value = 708.889453506068
$ws.Range("S7").Value = 0.02187326217661522
$ws.Range("T7").Value = 0.02187326217661522
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 29.80827733333334
$ws.Range("N8").Value = 89.42483200000001
$ws.Range("O8").Value = 0.4866975737940222
$ws.Range("P8").Value = 0.4866975737940221
$ws.Range("Q8").Value = 174.0799819273387
$ws.Range("R8").Value = 1566.719837346048
$ws.Range("S8").Value = 0.04834219720731215
$ws.Range("T8").Value = 0.04834219720731214
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.437868100938039
$ws.Range("P9").Value = 0.437868100938039
$ws.Range("S9").Value = 0.04349211343160785
$ws.Range("T9").Value = 0.04349211343160784
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("M10").Value = 4.620050333333332
$ws.Range("O10").Value = 0.07543432526793886
$ws.Range("P10").Value = 0.07543432526793886
$ws.Range("S10").Value = 0.007492663256724088
$ws.Range("T10").Value = 0.007492663256724087
